# Auto-generated edit script: refresh market-data derived columns (H-N)
# across the 8 job sheets to match the scheduled runner's latest pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 374.875  # H9: 414 -> 374.875
$ws.Cells.Item(9, 9).Value = 299.83334  # I9: 339.6 -> 299.83334
$ws.Cells.Item(9, 11).Value = 299.83334  # K9: 339.6 -> 299.83334
$ws.Cells.Item(9, 13).Value = -130.83334  # M9: -170.6 -> -130.83334
$ws.Cells.Item(64, 8).Value = 13666.667  # H64: 10601.111 -> 13666.667
$ws.Cells.Item(64, 9).Value = 12000  # I64: 6352.5 -> 12000
$ws.Cells.Item(64, 11).Value = 12000  # K64: 6352.5 -> 12000
$ws.Cells.Item(64, 13).Value = -11752  # M64: -6104.5 -> -11752
$ws.Cells.Item(67, 8).Value = 13666.667  # H67: 10601.111 -> 13666.667
$ws.Cells.Item(67, 9).Value = 12000  # I67: 6352.5 -> 12000
$ws.Cells.Item(67, 11).Value = 12000  # K67: 6352.5 -> 12000
$ws.Cells.Item(67, 13).Value = -11142  # M67: -5494.5 -> -11142
$ws.Cells.Item(75, 8).Value = 0  # H75: 54000 -> 0
$ws.Cells.Item(75, 10).Value = 0  # J75: 54000 -> 0
$ws.Cells.Item(75, 12).Value = 0  # L75: 54000 -> 0
$ws.Cells.Item(75, 14).Value = ""  # N75: clear (was -55872)
$ws.Cells.Item(78, 8).Value = 0  # H78: 54000 -> 0
$ws.Cells.Item(78, 10).Value = 0  # J78: 54000 -> 0
$ws.Cells.Item(78, 12).Value = 0  # L78: 162000 -> 0
$ws.Cells.Item(78, 14).Value = ""  # N78: clear (was -171360)
$ws.Cells.Item(96, 8).Value = 347.75  # H96: 348.08334 -> 347.75
$ws.Cells.Item(96, 9).Value = 403.9  # I96: 373.45456 -> 403.9
$ws.Cells.Item(96, 10).Value = 67  # J96: 69 -> 67
$ws.Cells.Item(96, 11).Value = 1211.7  # K96: 1120.36368 -> 1211.7
$ws.Cells.Item(96, 12).Value = 201  # L96: 207 -> 201
$ws.Cells.Item(96, 13).Value = 161.3000000000002  # M96: 252.6363200000001 -> 161.3000000000002
$ws.Cells.Item(96, 14).Value = -2947  # N96: -2953 -> -2947
$ws.Cells.Item(101, 8).Value = 3487.5386  # H101: 3197.3635 -> 3487.5386
$ws.Cells.Item(101, 9).Value = 593.55554  # I101: 646.75 -> 593.55554
$ws.Cells.Item(101, 11).Value = 1780.66662  # K101: 1940.25 -> 1780.66662
$ws.Cells.Item(101, 13).Value = -158.66662  # M101: -318.25 -> -158.66662
$ws.Cells.Item(112, 8).Value = 3556.0667  # H112: 3294.5557 -> 3556.0667
$ws.Cells.Item(112, 10).Value = 3797.6924  # J112: 3458.1875 -> 3797.6924
$ws.Cells.Item(112, 12).Value = 11393.0772  # L112: 10374.5625 -> 11393.0772
$ws.Cells.Item(112, 14).Value = -13609.0772  # N112: -12590.5625 -> -13609.0772
$ws.Cells.Item(138, 8).Value = 2612.1482  # H138: 2666.6667 -> 2612.1482
$ws.Cells.Item(138, 9).Value = 3154.3  # I138: 3283.889 -> 3154.3
$ws.Cells.Item(138, 10).Value = 2293.2354  # J138: 2296.3333 -> 2293.2354
$ws.Cells.Item(138, 11).Value = 9462.900000000001  # K138: 9851.667000000001 -> 9462.900000000001
$ws.Cells.Item(138, 12).Value = 6879.706200000001  # L138: 6888.999899999999 -> 6879.706200000001
$ws.Cells.Item(138, 13).Value = -4322.900000000001  # M138: -4711.667000000001 -> -4322.900000000001
$ws.Cells.Item(138, 14).Value = -17159.7062  # N138: -17168.9999 -> -17159.7062

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5035.4204  # H32: 4975.943 -> 5035.4204
$ws.Cells.Item(32, 9).Value = 5028.385  # I32: 4965.409 -> 5028.385
$ws.Cells.Item(32, 11).Value = 5028.385  # K32: 4965.409 -> 5028.385
$ws.Cells.Item(32, 13).Value = -4741.385  # M32: -4678.409 -> -4741.385
$ws.Cells.Item(101, 8).Value = 50000  # H101: 0 -> 50000
$ws.Cells.Item(101, 10).Value = 50000  # J101: 0 -> 50000
$ws.Cells.Item(101, 12).Value = 50000  # L101: 0 -> 50000
$ws.Cells.Item(101, 14).Value = -56490  # N101: None -> -56490
$ws.Cells.Item(133, 8).Value = 167666.67  # H133: 146999.5 -> 167666.67
$ws.Cells.Item(133, 10).Value = 167666.67  # J133: 146999.5 -> 167666.67
$ws.Cells.Item(133, 12).Value = 167666.67  # L133: 146999.5 -> 167666.67
$ws.Cells.Item(133, 14).Value = -172726.67  # N133: -152059.5 -> -172726.67

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4968.55  # H20: 5298.647 -> 4968.55
$ws.Cells.Item(20, 9).Value = 4091.4666  # I20: 4339.8335 -> 4091.4666
$ws.Cells.Item(20, 11).Value = 4091.4666  # K20: 4339.8335 -> 4091.4666
$ws.Cells.Item(20, 13).Value = -3844.4666  # M20: -4092.8335 -> -3844.4666
$ws.Cells.Item(76, 8).Value = 24166.5  # H76: 33333 -> 24166.5
$ws.Cells.Item(76, 10).Value = 24166.5  # J76: 33333 -> 24166.5
$ws.Cells.Item(76, 12).Value = 24166.5  # L76: 33333 -> 24166.5
$ws.Cells.Item(76, 14).Value = -24796.5  # N76: -33963 -> -24796.5
$ws.Cells.Item(79, 8).Value = 24166.5  # H79: 33333 -> 24166.5
$ws.Cells.Item(79, 10).Value = 24166.5  # J79: 33333 -> 24166.5
$ws.Cells.Item(79, 12).Value = 24166.5  # L79: 33333 -> 24166.5
$ws.Cells.Item(79, 14).Value = -26350.5  # N79: -35517 -> -26350.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 111.85714  # H7: 104.125 -> 111.85714
$ws.Cells.Item(7, 9).Value = 116.8  # I7: 105.666664 -> 116.8
$ws.Cells.Item(7, 11).Value = 116.8  # K7: 105.666664 -> 116.8
$ws.Cells.Item(7, 13).Value = -3.799999999999997  # M7: 7.333336000000003 -> -3.799999999999997
$ws.Cells.Item(22, 8).Value = 594.1667  # H22: 547.2857 -> 594.1667
$ws.Cells.Item(22, 9).Value = 258  # I22: 283 -> 258
$ws.Cells.Item(22, 10).Value = 762.25  # J22: 899.6667 -> 762.25
$ws.Cells.Item(22, 11).Value = 258  # K22: 283 -> 258
$ws.Cells.Item(22, 12).Value = 762.25  # L22: 899.6667 -> 762.25
$ws.Cells.Item(22, 13).Value = 92  # M22: 67 -> 92
$ws.Cells.Item(22, 14).Value = -1462.25  # N22: -1599.6667 -> -1462.25
$ws.Cells.Item(31, 8).Value = 4156.2856  # H31: 4237.353 -> 4156.2856
$ws.Cells.Item(31, 9).Value = 4079.4666  # I31: 4270.857 -> 4079.4666
$ws.Cells.Item(31, 11).Value = 4079.4666  # K31: 4270.857 -> 4079.4666
$ws.Cells.Item(31, 13).Value = -3784.4666  # M31: -3975.857 -> -3784.4666
$ws.Cells.Item(34, 8).Value = 4156.2856  # H34: 4237.353 -> 4156.2856
$ws.Cells.Item(34, 9).Value = 4079.4666  # I34: 4270.857 -> 4079.4666
$ws.Cells.Item(34, 11).Value = 4079.4666  # K34: 4270.857 -> 4079.4666
$ws.Cells.Item(34, 13).Value = -3877.4666  # M34: -4068.857 -> -3877.4666
$ws.Cells.Item(122, 8).Value = 4464.0713  # H122: 4245.933 -> 4464.0713
$ws.Cells.Item(122, 9).Value = 4388.6665  # I122: 4069 -> 4388.6665
$ws.Cells.Item(122, 11).Value = 13165.9995  # K122: 12207 -> 13165.9995
$ws.Cells.Item(122, 13).Value = -10715.9995  # M122: -9757 -> -10715.9995
$ws.Cells.Item(133, 8).Value = 58729.4  # H133: 63161.75 -> 58729.4
$ws.Cells.Item(133, 9).Value = 41000  # I133: 0 -> 41000
$ws.Cells.Item(133, 11).Value = 41000  # K133: 0 -> 41000
$ws.Cells.Item(133, 13).Value = -38470  # M133: None -> -38470
$ws.Cells.Item(141, 8).Value = 34480.75  # H141: 31836.4 -> 34480.75
$ws.Cells.Item(141, 10).Value = 39110.4  # J141: 39444 -> 39110.4
$ws.Cells.Item(141, 12).Value = 39110.4  # L141: 39444 -> 39110.4
$ws.Cells.Item(141, 14).Value = -49470.4  # N141: -49804 -> -49470.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 666966.7  # H7: 3000224.8 -> 666966.7
$ws.Cells.Item(7, 9).Value = 666966.7  # I7: 3000224.8 -> 666966.7
$ws.Cells.Item(7, 11).Value = 2000900.1  # K7: 9000674.399999999 -> 2000900.1
$ws.Cells.Item(7, 13).Value = -2000788.1  # M7: -9000562.399999999 -> -2000788.1
$ws.Cells.Item(11, 8).Value = 8334013  # H11: 8334021.5 -> 8334013
$ws.Cells.Item(11, 9).Value = 10000631  # I11: 11111734 -> 10000631
$ws.Cells.Item(11, 10).Value = 925  # J11: 884 -> 925
$ws.Cells.Item(11, 11).Value = 30001893  # K11: 33335202 -> 30001893
$ws.Cells.Item(11, 12).Value = 2775  # L11: 2652 -> 2775
$ws.Cells.Item(11, 13).Value = -30001753  # M11: -33335062 -> -30001753
$ws.Cells.Item(11, 14).Value = -3055  # N11: -2932 -> -3055
$ws.Cells.Item(12, 8).Value = 935.3684  # H12: 985.1111 -> 935.3684
$ws.Cells.Item(12, 10).Value = 992.4706  # J12: 1052 -> 992.4706
$ws.Cells.Item(12, 12).Value = 2977.4118  # L12: 3156 -> 2977.4118
$ws.Cells.Item(12, 14).Value = -3323.4118  # N12: -3502 -> -3323.4118
$ws.Cells.Item(130, 8).Value = 15000  # H130: 2000 -> 15000
$ws.Cells.Item(130, 10).Value = 15000  # J130: 2000 -> 15000
$ws.Cells.Item(130, 12).Value = 45000  # L130: 6000 -> 45000
$ws.Cells.Item(130, 14).Value = -55040  # N130: -16040 -> -55040
$ws.Cells.Item(140, 8).Value = 1728.3846  # H140: 1773.52 -> 1728.3846
$ws.Cells.Item(140, 9).Value = 1114.9231  # I140: 1157.8334 -> 1114.9231
$ws.Cells.Item(140, 11).Value = 3344.7693  # K140: 3473.5002 -> 3344.7693
$ws.Cells.Item(140, 13).Value = 1835.2307  # M140: 1706.4998 -> 1835.2307

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2843.3157  # H80: 2945.8125 -> 2843.3157
$ws.Cells.Item(80, 9).Value = 3023.5  # I80: 3144.6667 -> 3023.5
$ws.Cells.Item(80, 10).Value = 2338.8  # J80: 2349.25 -> 2338.8
$ws.Cells.Item(80, 11).Value = 3023.5  # K80: 3144.6667 -> 3023.5
$ws.Cells.Item(80, 12).Value = 2338.8  # L80: 2349.25 -> 2338.8
$ws.Cells.Item(80, 13).Value = -2025.5  # M80: -2146.6667 -> -2025.5
$ws.Cells.Item(80, 14).Value = -4334.8  # N80: -4345.25 -> -4334.8
$ws.Cells.Item(83, 8).Value = 2843.3157  # H83: 2945.8125 -> 2843.3157
$ws.Cells.Item(83, 9).Value = 3023.5  # I83: 3144.6667 -> 3023.5
$ws.Cells.Item(83, 10).Value = 2338.8  # J83: 2349.25 -> 2338.8
$ws.Cells.Item(83, 11).Value = 15117.5  # K83: 15723.3335 -> 15117.5
$ws.Cells.Item(83, 12).Value = 11694  # L83: 11746.25 -> 11694
$ws.Cells.Item(83, 13).Value = -10125.5  # M83: -10731.3335 -> -10125.5
$ws.Cells.Item(83, 14).Value = -21678  # N83: -21730.25 -> -21678
$ws.Cells.Item(113, 8).Value = 2766.3333  # H113: 2782.8333 -> 2766.3333
$ws.Cells.Item(113, 9).Value = 3449.5  # I113: 3474.25 -> 3449.5
$ws.Cells.Item(113, 11).Value = 3449.5  # K113: 3474.25 -> 3449.5
$ws.Cells.Item(113, 13).Value = -1279.5  # M113: -1304.25 -> -1279.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3332.3333  # H7: 3749.5 -> 3332.3333
$ws.Cells.Item(7, 9).Value = 2498.75  # I7: 2499.5 -> 2498.75
$ws.Cells.Item(7, 11).Value = 2498.75  # K7: 2499.5 -> 2498.75
$ws.Cells.Item(7, 13).Value = -2386.75  # M7: -2387.5 -> -2386.75
$ws.Cells.Item(16, 8).Value = 1980.6428  # H16: 2012.7858 -> 1980.6428
$ws.Cells.Item(16, 9).Value = 1659.3334  # I16: 1709.3334 -> 1659.3334
$ws.Cells.Item(16, 11).Value = 1659.3334  # K16: 1709.3334 -> 1659.3334
$ws.Cells.Item(16, 13).Value = -1489.3334  # M16: -1539.3334 -> -1489.3334
$ws.Cells.Item(22, 8).Value = 4172.3  # H22: 3911.182 -> 4172.3
$ws.Cells.Item(22, 9).Value = 2715.5  # I22: 3007 -> 2715.5
$ws.Cells.Item(22, 10).Value = 9999.5  # J22: 5493.5 -> 9999.5
$ws.Cells.Item(22, 11).Value = 2715.5  # K22: 3007 -> 2715.5
$ws.Cells.Item(22, 12).Value = 9999.5  # L22: 5493.5 -> 9999.5
$ws.Cells.Item(22, 13).Value = -2420.5  # M22: -2712 -> -2420.5
$ws.Cells.Item(22, 14).Value = -10589.5  # N22: -6083.5 -> -10589.5
$ws.Cells.Item(27, 8).Value = 4172.3  # H27: 3911.182 -> 4172.3
$ws.Cells.Item(27, 9).Value = 2715.5  # I27: 3007 -> 2715.5
$ws.Cells.Item(27, 10).Value = 9999.5  # J27: 5493.5 -> 9999.5
$ws.Cells.Item(27, 11).Value = 2715.5  # K27: 3007 -> 2715.5
$ws.Cells.Item(27, 12).Value = 9999.5  # L27: 5493.5 -> 9999.5
$ws.Cells.Item(27, 13).Value = -2608.5  # M27: -2900 -> -2608.5
$ws.Cells.Item(27, 14).Value = -10213.5  # N27: -5707.5 -> -10213.5
$ws.Cells.Item(38, 8).Value = 35000  # H38: 0 -> 35000
$ws.Cells.Item(38, 10).Value = 35000  # J38: 0 -> 35000
$ws.Cells.Item(38, 12).Value = 35000  # L38: 0 -> 35000
$ws.Cells.Item(38, 14).Value = -35820  # N38: None -> -35820
$ws.Cells.Item(40, 8).Value = 3138.8462  # H40: 3300.5 -> 3138.8462
$ws.Cells.Item(40, 9).Value = 2328.7144  # I40: 2517 -> 2328.7144
$ws.Cells.Item(40, 11).Value = 2328.7144  # K40: 2517 -> 2328.7144
$ws.Cells.Item(40, 13).Value = -2192.7144  # M40: -2381 -> -2192.7144
$ws.Cells.Item(61, 8).Value = 203605.2  # H61: 253839.75 -> 203605.2
$ws.Cells.Item(61, 10).Value = 2667  # J61: 0 -> 2667
$ws.Cells.Item(61, 12).Value = 2667  # L61: 0 -> 2667
$ws.Cells.Item(61, 14).Value = -3071  # N61: None -> -3071
$ws.Cells.Item(109, 8).Value = 40999.5  # H109: 41000 -> 40999.5
$ws.Cells.Item(109, 10).Value = 40999.5  # J109: 41000 -> 40999.5
$ws.Cells.Item(109, 12).Value = 40999.5  # L109: 41000 -> 40999.5
$ws.Cells.Item(109, 14).Value = -43773.5  # N109: -43774 -> -43773.5
$ws.Cells.Item(113, 8).Value = 203605.2  # H113: 253839.75 -> 203605.2
$ws.Cells.Item(113, 10).Value = 2667  # J113: 0 -> 2667
$ws.Cells.Item(113, 12).Value = 2667  # L113: 0 -> 2667
$ws.Cells.Item(113, 14).Value = -7007  # N113: None -> -7007
$ws.Cells.Item(122, 8).Value = 3232.3914  # H122: 3304.75 -> 3232.3914
$ws.Cells.Item(122, 9).Value = 3216.2222  # I122: 3274.5 -> 3216.2222
$ws.Cells.Item(122, 10).Value = 3290.6  # J122: 3425.75 -> 3290.6
$ws.Cells.Item(122, 11).Value = 9648.6666  # K122: 9823.5 -> 9648.6666
$ws.Cells.Item(122, 12).Value = 9871.799999999999  # L122: 10277.25 -> 9871.799999999999
$ws.Cells.Item(122, 13).Value = -7198.6666  # M122: -7373.5 -> -7198.6666
$ws.Cells.Item(122, 14).Value = -14771.8  # N122: -15177.25 -> -14771.8
$ws.Cells.Item(126, 8).Value = 3332.3333  # H126: 3749.5 -> 3332.3333
$ws.Cells.Item(126, 9).Value = 2498.75  # I126: 2499.5 -> 2498.75
$ws.Cells.Item(126, 11).Value = 7496.25  # K126: 7498.5 -> 7496.25
$ws.Cells.Item(126, 13).Value = -5026.25  # M126: -5028.5 -> -5026.25
$ws.Cells.Item(136, 8).Value = 4514.7407  # H136: 4430.4287 -> 4514.7407
$ws.Cells.Item(136, 9).Value = 3852.3914  # I136: 3781.625 -> 3852.3914
$ws.Cells.Item(136, 11).Value = 11557.1742  # K136: 11344.875 -> 11557.1742
$ws.Cells.Item(136, 13).Value = -9007.174199999999  # M136: -8794.875 -> -9007.174199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 877  # H107: 648.1667 -> 877
$ws.Cells.Item(107, 9).Value = 773.3333  # I107: 609.5 -> 773.3333
$ws.Cells.Item(107, 10).Value = 1499  # J107: 841.5 -> 1499
$ws.Cells.Item(107, 11).Value = 2319.9999  # K107: 1828.5 -> 2319.9999
$ws.Cells.Item(107, 12).Value = 4497  # L107: 2524.5 -> 4497
$ws.Cells.Item(107, 13).Value = -399.9998999999998  # M107: 91.5 -> -399.9998999999998
$ws.Cells.Item(107, 14).Value = -8337  # N107: -6364.5 -> -8337
$ws.Cells.Item(122, 8).Value = 9028.714  # H122: 10599.6 -> 9028.714
$ws.Cells.Item(122, 9).Value = 8400.666999999999  # I122: 10499 -> 8400.666999999999
$ws.Cells.Item(122, 10).Value = 9499.75  # J122: 10666.667 -> 9499.75
$ws.Cells.Item(122, 11).Value = 25202.001  # K122: 31497 -> 25202.001
$ws.Cells.Item(122, 12).Value = 28499.25  # L122: 32000.001 -> 28499.25
$ws.Cells.Item(122, 13).Value = -22752.001  # M122: -29047 -> -22752.001
$ws.Cells.Item(122, 14).Value = -33399.25  # N122: -36900.001 -> -33399.25
$ws.Cells.Item(126, 8).Value = 3450  # H126: 3159.6 -> 3450
$ws.Cells.Item(126, 9).Value = 3450  # I126: 3199.5 -> 3450
$ws.Cells.Item(126, 10).Value = 0  # J126: 3000 -> 0
$ws.Cells.Item(126, 11).Value = 10350  # K126: 9598.5 -> 10350
$ws.Cells.Item(126, 12).Value = 0  # L126: 9000 -> 0
$ws.Cells.Item(126, 13).Value = -7880  # M126: -7128.5 -> -7880
$ws.Cells.Item(126, 14).Value = ""  # N126: clear (was -13940)
